# Update "想去人数" (F column) counts across sheets, as generated by the
# gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7743
$ws1.Range("F5").Value = 8070
$ws1.Range("F7").Value = 592
$ws1.Range("F8").Value = 631
$ws1.Range("F10").Value = 137
$ws1.Range("F13").Value = 36
$ws1.Range("F17").Value = 263
$ws1.Range("F20").Value = 147
$ws1.Range("F22").Value = 77
$ws1.Range("F23").Value = 614
$ws1.Range("F24").Value = 2202
$ws1.Range("F25").Value = 734
$ws1.Range("F29").Value = 612

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 323

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 449

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 449
$ws4.Range("F3").Value = 7743
$ws4.Range("F7").Value = 8074
$ws4.Range("F9").Value = 592
$ws4.Range("F10").Value = 631
$ws4.Range("F13").Value = 137
$ws4.Range("F15").Value = 323
$ws4.Range("F19").Value = 36
$ws4.Range("F26").Value = 263
$ws4.Range("F29").Value = 147
$ws4.Range("F31").Value = 77
$ws4.Range("F32").Value = 614
$ws4.Range("F33").Value = 2202
$ws4.Range("F34").Value = 734
$ws4.Range("F39").Value = 612

$wb.Save()
